# Apply cell value updates from the cryptos data refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.913.00"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "3.163.43"
$ws.Range("E3").Value = "  +4.19%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'579.05"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("D6").Value = "'150.47"
$ws.Range("E6").Value = "  +7.01%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.163.08"
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  +6.69%  "
$ws.Range("D11").Value = "'6.20"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "'0.500"
$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  +18.90%  "
$ws.Range("D14").Value = "'37.48"
$ws.Range("E14").Value = "  +6.18%  "
$ws.Range("D15").Value = "3.681.06"
$ws.Range("E15").Value = "  +4.20%  "
$ws.Range("D16").Value = "64.982.96"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "3.171.25"
$ws.Range("E17").Value = "  +4.62%  "
$ws.Range("E18").Value = "  +6.31%  "
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").Value = "'510.05"
$ws.Range("E20").Value = "  +8.11%  "
$ws.Range("D21").Value = "'14.83"
$ws.Range("E21").Value = "  +6.10%  "
$ws.Range("E22").Value = "  +6.88%  "
$ws.Range("D23").Value = "'15.31"
$ws.Range("E23").Value = "  +5.84%  "
$ws.Range("E24").Value = "  +4.03%  "
$ws.Range("D25").Value = "'84.99"
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'9.05"
$ws.Range("E27").Value = "  +12.53%  "
$ws.Range("E28").Value = "  +5.29%  "
$ws.Range("E29").Value = "  +8.18%  "
$ws.Range("E30").Value = "  +6.81%  "
$ws.Range("D31").Value = "'2.78"
$ws.Range("E31").Value = "  +15.28%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'1.20"
$ws.Range("E33").Value = "  +4.91%  "
$ws.Range("D34").Value = "'6.30"
$ws.Range("E34").Value = "  +11.75%  "
$ws.Range("E35").Value = "  +6.79%  "
$ws.Range("D36").Value = "'55.71"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("E37").Value = "  +10.58%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'471.91"
$ws.Range("E38").Value = "  +7.62%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.13"
$ws.Range("E39").Value = "  +14.19%  "
$ws.Range("D40").Value = "'0.0420"
$ws.Range("E40").Value = "  +3.72%  "
$ws.Range("E41").Value = "  +4.79%  "
$ws.Range("D42").Value = "3.063.89"
$ws.Range("E42").Value = "  +2.39%  "
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("E44").Value = "  +5.79%  "
$ws.Range("E45").Value = "  +8.38%  "
$ws.Range("D46").Value = "'29.31"
$ws.Range("E46").Value = "  +6.68%  "
$ws.Range("D47").Value = "0.0₃0601"
$ws.Range("E47").Value = "  +18.43%  "
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "'0.115"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("E50").Value = "  +8.96%  "
$ws.Range("D51").Value = "'119.73"
$ws.Range("E51").Value = "  +1.41%  "
